# Updating state policy data - new states issued closure orders
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Florida (row 11): statewide complete closure order issued -> add complete closure date,
# remove the now-obsolete note about no statewide order
$ws.Range("D11").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value = 43924
$ws.Range("G11").ClearContents()

# Georgia (row 12): complete closure date added
$ws.Range("D12").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value = 43924

# South Carolina (row 42): complete closure date added
$ws.Range("D42").Copy()
$ws.Range("F42").PasteSpecial(-4122)
$ws.Range("F42").Value = 43921

# Texas (row 45): complete closure date added
$ws.Range("D45").Copy()
$ws.Range("F45").PasteSpecial(-4122)
$ws.Range("F45").Value = 43923

# Reflect final cursor/selection position used while entering this data
$null = $ws.Range("F52").Select()
